$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "rank" column (E) previously stored JSON scoring rules using a
# "regex" key. The rank is now inserted as a plain string in the DB, so
# the stored JSON rule objects use a "query" key instead of "regex".
# Update every data row (2-10) in column E accordingly.

$ws.Range("E2").Value = '[{"query":{"status":{"$ne":"reachable"}},"score":0},{"query":{"status":"reachable"},"score":100}]'
$ws.Range("E3").Value = '[{"query":{"interface":"down"},"score":0},{"query":{"interface":"up"},"score":100}]'
$ws.Range("E4").Value = '[{"query":{"interface":"down"},"score":0},{"query":{"interface":"up"},"score":100}]'
$ws.Range("E5").Value = '[{"query":{"interface":"down"},"score":0},{"query":{"interface":"up"},"score":100}]'
$ws.Range("E6").Value = '[{"query":{"status":"down"},"score":0},{"query":{"status":"reachable"},"score":100}]'
$ws.Range("E7").Value = '[{"query":{"interface":"down"},"score":0},{"query":{"interface":"up"},"score":100}]'
$ws.Range("E8").Value = '[{"query":{"interface":"down"},"score":0},{"query":{"interface":"up"},"score":100}]'
$ws.Range("E9").Value = '[{"query":{"interface":"down"},"score":0},{"query":{"interface":"up"},"score":100}]'
$ws.Range("E10").Value = '[{"query":{"result":"failed"},"score":0},{"query":{"result":"success"},"score":100}]'

# Reflect the last-edited cell selection recorded in the workbook.
$ws.Range("E8").Select()
